# Auto-generated edit script applying scheduled runner updates to Jenova_Profits sheets
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 20203
$ws.Range("I69").Value = 16500
$ws.Range("J69").Value = 22671.666
$ws.Range("K69").Value = 49500
$ws.Range("L69").Value = 68014.99800000001
$ws.Range("M69").Value = -48626
$ws.Range("N69").Value = -69762.99800000001
$ws.Range("H70").Value = 144250.58
$ws.Range("I70").Value = 350
$ws.Range("J70").Value = 168234
$ws.Range("K70").Value = 1050
$ws.Range("L70").Value = 504702
$ws.Range("M70").Value = -780
$ws.Range("N70").Value = -505242
$ws.Range("H72").Value = 20203
$ws.Range("I72").Value = 16500
$ws.Range("J72").Value = 22671.666
$ws.Range("K72").Value = 148500
$ws.Range("L72").Value = 204044.994
$ws.Range("M72").Value = -144132
$ws.Range("N72").Value = -212780.994
$ws.Range("H73").Value = 144250.58
$ws.Range("I73").Value = 350
$ws.Range("J73").Value = 168234
$ws.Range("K73").Value = 1050
$ws.Range("L73").Value = 504702
$ws.Range("M73").Value = -114
$ws.Range("N73").Value = -506574
$ws.Range("H74").Value = 8061.385
$ws.Range("I74").Value = 6572.5454
$ws.Range("K74").Value = 6572.5454
$ws.Range("M74").Value = -5636.5454
$ws.Range("H77").Value = 8061.385
$ws.Range("I77").Value = 6572.5454
$ws.Range("K77").Value = 32862.727
$ws.Range("M77").Value = -28182.727
$ws.Range("H100").Value = 2220.3333
$ws.Range("I100").Value = 2197.6155
$ws.Range("J100").Value = 2279.4
$ws.Range("K100").Value = 2197.6155
$ws.Range("L100").Value = 2279.4
$ws.Range("M100").Value = -1656.6155
$ws.Range("N100").Value = -3361.4
$ws.Range("H106").Value = 2890.7273
$ws.Range("I106").Value = 3399.7144
$ws.Range("K106").Value = 3399.7144
$ws.Range("M106").Value = -2768.7144
$ws.Range("H125").Value = 7233.091
$ws.Range("I125").Value = 6803.5
$ws.Range("J125").Value = 8378.666999999999
$ws.Range("K125").Value = 61231.5
$ws.Range("L125").Value = 75408.003
$ws.Range("M125").Value = -58771.5
$ws.Range("N125").Value = -80328.003
$ws.Range("H131").Value = 4215.4
$ws.Range("I131").Value = 2644.2
$ws.Range("K131").Value = 7932.599999999999
$ws.Range("M131").Value = -2892.599999999999
$ws.Range("H132").Value = 2994.5417
$ws.Range("I132").Value = 3055.127
$ws.Range("J132").Value = 2570.4443
$ws.Range("K132").Value = 9165.380999999999
$ws.Range("L132").Value = 7711.3329
$ws.Range("M132").Value = -6635.380999999999
$ws.Range("N132").Value = -12771.3329
$ws.Range("H137").Value = 6107.316
$ws.Range("I137").Value = 7950
$ws.Range("J137").Value = 5256.846
$ws.Range("K137").Value = 23850
$ws.Range("L137").Value = 15770.538
$ws.Range("M137").Value = -21300
$ws.Range("N137").Value = -20870.538
$ws.Range("H138").Value = 7366.0835
$ws.Range("I138").Value = 6002.364
$ws.Range("J138").Value = 8520
$ws.Range("K138").Value = 18007.092
$ws.Range("L138").Value = 25560
$ws.Range("M138").Value = -12867.092
$ws.Range("N138").Value = -35840

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3093.377
$ws.Range("I32").Value = 2391.3276
$ws.Range("K32").Value = 2391.3276
$ws.Range("M32").Value = -2104.3276
$ws.Range("H47").Value = 50000
$ws.Range("I47").Value = 50000
$ws.Range("K47").Value = 50000
$ws.Range("M47").Value = -49275
$ws.Range("H61").Value = 4280.136
$ws.Range("I61").Value = 2051.077
$ws.Range("J61").Value = 7499.8887
$ws.Range("K61").Value = 2051.077
$ws.Range("L61").Value = 7499.8887
$ws.Range("M61").Value = -1839.077
$ws.Range("N61").Value = -7923.8887
$ws.Range("H102").Value = 2126.4
$ws.Range("I102").Value = 2409.1428
$ws.Range("J102").Value = 1466.6666
$ws.Range("K102").Value = 2409.1428
$ws.Range("L102").Value = 1466.6666
$ws.Range("M102").Value = -787.1428000000001
$ws.Range("N102").Value = -4710.6666
$ws.Range("H132").Value = 4887.4385
$ws.Range("I132").Value = 1529.3414
$ws.Range("J132").Value = 13492.5625
$ws.Range("K132").Value = 4588.0242
$ws.Range("L132").Value = 40477.6875
$ws.Range("M132").Value = -2058.0242
$ws.Range("N132").Value = -45537.6875
$ws.Range("H136").Value = 4280.136
$ws.Range("I136").Value = 2051.077
$ws.Range("J136").Value = 7499.8887
$ws.Range("K136").Value = 6153.231000000001
$ws.Range("L136").Value = 22499.6661
$ws.Range("M136").Value = -3603.231000000001
$ws.Range("N136").Value = -27599.6661

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 23823.756
$ws.Range("I134").Value = 2940.4595
$ws.Range("J134").Value = 88213.914
$ws.Range("K134").Value = 8821.378499999999
$ws.Range("L134").Value = 264641.742
$ws.Range("M134").Value = -6286.378499999999
$ws.Range("N134").Value = -269711.742

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4853.7646
$ws.Range("I31").Value = 1720.5714
$ws.Range("J31").Value = 7047
$ws.Range("K31").Value = 1720.5714
$ws.Range("L31").Value = 7047
$ws.Range("M31").Value = -1425.5714
$ws.Range("N31").Value = -7637
$ws.Range("H34").Value = 4853.7646
$ws.Range("I34").Value = 1720.5714
$ws.Range("J34").Value = 7047
$ws.Range("K34").Value = 1720.5714
$ws.Range("L34").Value = 7047
$ws.Range("M34").Value = -1518.5714
$ws.Range("N34").Value = -7451
$ws.Range("H105").Value = 2675.6667
$ws.Range("I105").Value = 2675.6667
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 2675.6667
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -928.6667000000002
$ws.Range("N105").ClearContents()
$ws.Range("H132").Value = 3490.3704
$ws.Range("I132").Value = 3811.238
$ws.Range("J132").Value = 2367.3333
$ws.Range("K132").Value = 11433.714
$ws.Range("L132").Value = 7101.999899999999
$ws.Range("M132").Value = -8903.714
$ws.Range("N132").Value = -12161.9999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 51451.562
$ws.Range("I5").Value = 134666.5
$ws.Range("J5").Value = 1522.6
$ws.Range("K5").Value = 403999.5
$ws.Range("L5").Value = 4567.799999999999
$ws.Range("M5").Value = -403887.5
$ws.Range("N5").Value = -4791.799999999999
$ws.Range("H68").Value = 2555.1428
$ws.Range("J68").Value = 3356.2
$ws.Range("L68").Value = 10068.6
$ws.Range("N68").Value = -11690.6
$ws.Range("H71").Value = 2555.1428
$ws.Range("J71").Value = 3356.2
$ws.Range("L71").Value = 30205.8
$ws.Range("N71").Value = -38317.8
$ws.Range("H107").Value = 5959
$ws.Range("I107").Value = 1019.75
$ws.Range("K107").Value = 3059.25
$ws.Range("M107").Value = -1139.25
$ws.Range("H132").Value = 3497.3333
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 3497.3333
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 31475.9997
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -36535.9997
$ws.Range("H135").Value = 51451.562
$ws.Range("I135").Value = 134666.5
$ws.Range("J135").Value = 1522.6
$ws.Range("K135").Value = 1211998.5
$ws.Range("L135").Value = 13703.4
$ws.Range("M135").Value = -1209463.5
$ws.Range("N135").Value = -18773.4
$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("M138").ClearContents()
$ws.Range("N138").ClearContents()
$ws.Range("H139").Value = 1616.6875
$ws.Range("I139").Value = 1616.6875
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 4850.0625
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = 289.9375
$ws.Range("N139").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H133").Value = 62000
$ws.Range("J133").Value = 62000
$ws.Range("L133").Value = 62000
$ws.Range("N133").Value = -72120

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 423556.66
$ws.Range("I7").Value = 839638.3
$ws.Range("K7").Value = 839638.3
$ws.Range("M7").Value = -839526.3
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("H46").Value = 3283.5676
$ws.Range("J46").Value = 3755.125
$ws.Range("L46").Value = 3755.125
$ws.Range("N46").Value = -4131.125
$ws.Range("H126").Value = 423556.66
$ws.Range("I126").Value = 839638.3
$ws.Range("K126").Value = 2518914.9
$ws.Range("M126").Value = -2516444.9
$ws.Range("H132").Value = 10949.75
$ws.Range("I132").Value = 7399.5
$ws.Range("J132").Value = 14500
$ws.Range("K132").Value = 22198.5
$ws.Range("L132").Value = 43500
$ws.Range("M132").Value = -19668.5
$ws.Range("N132").Value = -48560
$ws.Range("H136").Value = 4254
$ws.Range("I136").Value = 3248.5715
$ws.Range("J136").Value = 6600
$ws.Range("K136").Value = 9745.7145
$ws.Range("L136").Value = 19800
$ws.Range("M136").Value = -7195.7145
$ws.Range("N136").Value = -24900

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6614.846
$ws.Range("I62").Value = 6070
$ws.Range("J62").Value = 6955.375
$ws.Range("K62").Value = 6070
$ws.Range("L62").Value = 6955.375
$ws.Range("M62").Value = -5446
$ws.Range("N62").Value = -8203.375
$ws.Range("H65").Value = 6614.846
$ws.Range("I65").Value = 6070
$ws.Range("J65").Value = 6955.375
$ws.Range("K65").Value = 30350
$ws.Range("L65").Value = 34776.875
$ws.Range("M65").Value = -27230
$ws.Range("N65").Value = -41016.875
$ws.Range("H113").Value = 768.075
$ws.Range("I113").Value = 617.09375
$ws.Range("J113").Value = 1372
$ws.Range("K113").Value = 1851.28125
$ws.Range("L113").Value = 4116
$ws.Range("M113").Value = 318.71875
$ws.Range("N113").Value = -8456
$ws.Range("H132").Value = 24926.98
$ws.Range("I132").Value = 4316.5127
$ws.Range("J132").Value = 98000.45
$ws.Range("K132").Value = 12949.5381
$ws.Range("L132").Value = 294001.35
$ws.Range("M132").Value = -10419.5381
$ws.Range("N132").Value = -299061.35

